# Weekly update: insert a new price-report row for "Feria Lagunitas de
# Puerto Montt - Zanahoria" above the current row 191, pushing the
# existing historical rows (191-220) down by one (to 192-221).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record; Excel copies the formatting (incl. the
# date number format on column D) from the surrounding rows automatically.
$ws.Rows("191:191").Insert()

# Populate the newly inserted row with this week's data point.
$ws.Cells.Item(191, 1).Value = 4
$ws.Cells.Item(191, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(191, 3).Value = "Los Lagos"
$ws.Cells.Item(191, 4).Value = 44474
$ws.Cells.Item(191, 5).Value = 10
$ws.Cells.Item(191, 6).Value = 100114013
$ws.Cells.Item(191, 7).Value = "Zanahoria"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 850
$ws.Cells.Item(191, 11).Value = 12000
$ws.Cells.Item(191, 12).Value = 12000
$ws.Cells.Item(191, 13).Value = 12000
$ws.Cells.Item(191, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(191, 15).Value = "Región de Ñuble"
$ws.Cells.Item(191, 16).Value = 600
$ws.Cells.Item(191, 17).Value = 20
$ws.Cells.Item(191, 18).Value = "Hortaliza"
